$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 0.002
$ws.Range("K14").Value = 1007
$ws.Range("L14").Value = 0.002014
